# feature: handle \n in file excel
#
# The "No precondition" text in C2 is replaced by a two-line value
# ("Condition 1" / "Condition 2"); wrap text is turned on for that cell so
# both lines are visible, which also widens column C and grows row 2's
# height, and the active selection moves to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multi-line content (embedded newline -> wraps inside the cell).
$ws.Range("C2").Value = "Condition 1`nCondition 2"

# Show every line of the new text.
$ws.Range("C2").WrapText = $true

# Make room for the wrapped text.
$ws.Columns("C").ColumnWidth = 33.42578125
$ws.Rows("2").RowHeight = 30

# Cursor ends up on the edited cell.
$ws.Range("C2").Select()
